# Electives Replace in User Schedule
# - Use the case with the best match score to replace electives in the
#   Web and CyberSecurity schedule rows (append the newly-chosen elective
#   course to the existing comma-separated list).
# - Add (select) the track button cells next to the merged "E - elective"
#   header so the control has backing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Web track, Semester 8): best-match case adds CPSC 3118E elective.
$ws.Range("H6").Value = "CPSC 4125, CPSC 4175, CPSC 4115, GE, GE, CPSC 3118E"

# Row 7 (CyberSecurity track, Semester 7 & 8): best-match cases add
# CYBR 3136E and CYBR 3126E electives respectively.
$ws.Range("H7").Value = "CPSC 4157, CPSC 4127, CPSC 4130, CPSC 2125, CYBR 3136E"
$ws.Range("I7").Value = "CYBR 4160, CYBR 4166, CPSC 4000, CPSC 4138, CPSC 3118, CYBR 3126E"

# Give the track-button helper cells (inside the merged I2:K2 range) a
# real backing style so the new track button renders correctly.
$ws.Range("J2").Style = "Normal"
$ws.Range("K2").Style = "Normal"
